$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "114÷9=12, 6"
$t.Cell(1, 2).Range.Text = "965÷7=137, 6"
$t.Cell(1, 3).Range.Text = "135÷5=27, 0"
$t.Cell(1, 4).Range.Text = "231÷7=33, 0"
$t.Cell(1, 5).Range.Text = "325÷3=108, 1"

$t.Cell(5, 1).Range.Text = "855÷5=171, 0"
$t.Cell(5, 2).Range.Text = "470÷3=156, 2"
$t.Cell(5, 3).Range.Text = "397÷2=198, 1"
$t.Cell(5, 4).Range.Text = "661÷8=82, 5"
$t.Cell(5, 5).Range.Text = "102÷6=17, 0"

$t.Cell(9, 1).Range.Text = "106÷9=11, 7"
$t.Cell(9, 2).Range.Text = "988÷4=247, 0"
$t.Cell(9, 3).Range.Text = "681÷7=97, 2"
$t.Cell(9, 4).Range.Text = "196÷5=39, 1"
$t.Cell(9, 5).Range.Text = "478÷9=53, 1"

$t.Cell(13, 1).Range.Text = "177÷5=35, 2"
$t.Cell(13, 2).Range.Text = "686÷9=76, 2"
$t.Cell(13, 3).Range.Text = "592÷2=296, 0"
$t.Cell(13, 4).Range.Text = "626÷2=313, 0"
$t.Cell(13, 5).Range.Text = "288÷3=96, 0"

$t.Cell(17, 1).Range.Text = "373÷5=74, 3"
$t.Cell(17, 2).Range.Text = "899÷8=112, 3"
$t.Cell(17, 3).Range.Text = "271÷3=90, 1"
$t.Cell(17, 4).Range.Text = "955÷4=238, 3"
$t.Cell(17, 5).Range.Text = "987÷8=123, 3"

Write-Host "All cells updated"
